# Generate Report for Handback
#
# The handback-status job re-ran for 39c72983-74a4-4de4-95bb-f2d143388abb
# and found the handback file in the repo was stale (not the newest
# commit). This records that result on row 6 ("39c72983-...") of both the
# zh-cn and de-de status sheets:
#   - Latest Target File (I) gets a hyperlink to the (now-fetched) md file
#   - Latest Handback File (J) gets the generated xlf filename
#   - Latest Handback DateTime (K) gets the generation timestamp
#   - Error Detail (P) records the stale-version message
# Column P is also widened so the long error text is readable.

$wb = $excel.ActiveWorkbook

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9069659986ef6a95b4ab247d5e75d119b9bba7fe/e2e/39c72983-74a4-4de4-95bb-f2d143388abb.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3bc4d2aed7b226992f0fbb213dd057d5b4cfba0/e2e/39c72983-74a4-4de4-95bb-f2d143388abb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9069659986ef6a95b4ab247d5e75d119b9bba7fe/e2e/39c72983-74a4-4de4-95bb-f2d143388abb.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J6").Value = "39c72983-74a4-4de4-95bb-f2d143388abb.eaab5ab5fde9b99f06dead5efc5098284745664a.zh-cn.xlf"
$wsZh.Range("K6").Value = "2016-08-28 18:42:06"
$wsZh.Range("P6").Value = $errorDetail

$wsZh.Range("I6").Value = "39c72983-74a4-4de4-95bb-f2d143388abb.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I6"), $latestMdUrl, "", "", "39c72983-74a4-4de4-95bb-f2d143388abb.md")

$wsZh.Columns.Item(16).ColumnWidth = 40

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J6").Value = "39c72983-74a4-4de4-95bb-f2d143388abb.eaab5ab5fde9b99f06dead5efc5098284745664a.de-de.xlf"
$wsDe.Range("K6").Value = "2016-08-28 18:42:13"
$wsDe.Range("P6").Value = $errorDetail

$wsDe.Range("I6").Value = "39c72983-74a4-4de4-95bb-f2d143388abb.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I6"), $latestMdUrl, "", "", "39c72983-74a4-4de4-95bb-f2d143388abb.md")

$wsDe.Columns.Item(16).ColumnWidth = 40
